$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.6753301551942219
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 26.21740644021617
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 29.05731038884925

# Row 3
$ws.Range("B3").Value = 1.459612070389937
$ws.Range("C3").Value = 0.3127903958511391
$ws.Range("D3").Value = 26.21740644021617
$ws.Range("E3").Value = 645.3272768299601
$ws.Range("G3").Value = 673.3170857364173

# Row 4
$ws.Range("B4").Value = 3.230985683306322
$ws.Range("C4").Value = 1.667794583268128
$ws.Range("D4").Value = 0.8054896365839992
$ws.Range("E4").Value = 0.496779210170732
$ws.Range("G4").Value = 6.201049113329182

# Row 5
$ws.Range("B5").Value = 3.230985683306322
$ws.Range("C5").Value = 1.667794583268128
$ws.Range("D5").Value = 0.1575252929769615
$ws.Range("E5").Value = 8.660232485948974
$ws.Range("G5").Value = 13.71653804550039

# Row 6
$ws.Range("B6").Value = 3.230985683306322
$ws.Range("C6").Value = 1.667794583268128
$ws.Range("D6").Value = 26.21740644021617
$ws.Range("E6").Value = 0.496779210170732
$ws.Range("G6").Value = 31.61296591696135

# Row 7
$ws.Range("B7").Value = 0.127881588408715
$ws.Range("C7").Value = 0.00007097389502863649
$ws.Range("D7").Value = 0.8054896365839992
$ws.Range("E7").Value = 0.496779210170732
$ws.Range("G7").Value = 1.430221409058475
